$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 189, shifting existing rows 189:322 down to 190:323.
$ws.Rows.Item(189).EntireRow.Insert()

# Populate the newly inserted row 189 with the new data record.
$ws.Cells.Item(189, 1).Value = 9
$ws.Cells.Item(189, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(189, 3).Value = "Metropolitana"
$ws.Cells.Item(189, 4).Value = 44824
$ws.Cells.Item(189, 5).Value = 13
$ws.Cells.Item(189, 6).Value = 300000001
$ws.Cells.Item(189, 7).Value = "Rabanito"
$ws.Cells.Item(189, 8).Value = "Sin especificar"
$ws.Cells.Item(189, 9).Value = "Primera"
$ws.Cells.Item(189, 10).Value = 6100
$ws.Cells.Item(189, 11).Value = 3000
$ws.Cells.Item(189, 12).Value = 3000
$ws.Cells.Item(189, 13).Value = 3000
$ws.Cells.Item(189, 14).Value = "$/cien unidades (volumen en unidades)"
$ws.Cells.Item(189, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(189, 16).Value = 30
$ws.Cells.Item(189, 17).Value = 100
$ws.Cells.Item(189, 18).Value = "Hortaliza"
